$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 18

$ws.Cells.Item($row, 1).Value = 11
$ws.Cells.Item($row, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item($row, 3).Value = "Bíobío"

$ws.Cells.Item($row, 4).Value = 44505
$ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item($row, 5).Value = 8
$ws.Cells.Item($row, 6).Value = 100112026
$ws.Cells.Item($row, 7).Value = "Haba"
$ws.Cells.Item($row, 8).Value = "Sin especificar"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 180
$ws.Cells.Item($row, 11).Value = 6000
$ws.Cells.Item($row, 12).Value = 6500
$ws.Cells.Item($row, 13).Value = 6222
$ws.Cells.Item($row, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item($row, 15).Value = "Región del Maule"
$ws.Cells.Item($row, 16).Value = 249
$ws.Cells.Item($row, 17).Value = 25
$ws.Cells.Item($row, 18).Value = "Hortaliza"
